# Applies cryptocurrency price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.140.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.321.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "

# Row 7
$ws.Range("E7").Value = "  +0.28%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.28%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.11%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "

# Row 12
$ws.Range("E12").Value = "  -0.86%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.91%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.685.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.329.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.59%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.067.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.36%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.17%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0913"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.83%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.02%  "

# Row 24
$ws.Range("E24").Value = "  -1.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.91%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.41%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.26%  "

# Row 30
$ws.Range("E30").Value = "  +0.70%  "

# Row 31
$ws.Range("E31").Value = "  -1.88%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.76%  "

# Row 33
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.45%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0699"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.14%  "

# Row 38
$ws.Range("E38").Value = "  +2.00%  "

# Row 39
$ws.Range("E39").Value = "  +0.12%  "

# Row 40
$ws.Range("E40").Value = "  -0.33%  "

# Row 41
$ws.Range("E41").Value = "  +0.26%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.995.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.26%  "

# Row 43
$ws.Range("E43").Value = "  +1.73%  "

# Row 44
$ws.Range("E44").Value = "  -4.55%  "

# Row 45
$ws.Range("E45").Value = "  +1.05%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.55%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.38%  "

# Row 48
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.07%  "

# Row 49
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.70%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.549.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.95%  "

# Row 51
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.44%  "
